$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("book")

for ($r = 2; $r -le 137; $r++) {
    $ws.Cells.Item($r, 10).Value2 = 0
}
